# Recolor every slide's background fill from the old purple (7B5EA8)
# to the new blue (0075B2).
$p = $ppt.ActivePresentation

# VBA/PowerPoint RGB() is not available as a PowerShell cmdlet here, so the
# RGB(0x00, 0x75, 0xB2) value is pre-computed as R + G*256 + B*65536.
$newColor = 11695360

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $slide.Background.Fill.Solid()
    $slide.Background.Fill.ForeColor.RGB = $newColor
}

Write-Host "Updated background color on $($p.Slides.Count) slides"
